$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.728.72'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '3.153.62'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.51'
$ws.Range('E5').Value = '  +2.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.17'
$ws.Range('E6').Value = '  +5.71%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.149.86'
$ws.Range('E8').Value = '  +2.53%  '
$ws.Range('E9').Value = '  +4.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('E10').Value = '  +6.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.18'
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('E12').Value = '  +7.52%  '
$ws.Range('E13').Value = '  +12.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.91'
$ws.Range('E14').Value = '  +8.62%  '
$ws.Range('D15').Value = '3.670.36'
$ws.Range('E15').Value = '  +2.41%  '
$ws.Range('D16').Value = '64.855.61'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.23'
$ws.Range('E17').Value = '  +7.69%  '
$ws.Range('D18').Value = '3.156.07'
$ws.Range('E18').Value = '  +2.59%  '
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '518.05'
$ws.Range('E20').Value = '  +8.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.01'
$ws.Range('E21').Value = '  +8.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.738'
$ws.Range('E22').Value = '  +10.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.21'
$ws.Range('E23').Value = '  +8.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.84'
$ws.Range('E24').Value = '  +4.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.26'
$ws.Range('E25').Value = '  +5.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  +5.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.76'
$ws.Range('E28').Value = '  +10.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  +6.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '27.98'
$ws.Range('E30').Value = '  +6.62%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.66'
$ws.Range('E32').Value = '  +8.07%  '
$ws.Range('B33').Value = 'Mantle'
$ws.Range('C33').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.19'
$ws.Range('E33').Value = '  +4.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.13'
$ws.Range('E34').Value = '  +10.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.58'
$ws.Range('E35').Value = '  +6.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.70'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '487.44'
$ws.Range('E37').Value = '  +8.66%  '
$ws.Range('E38').Value = '  +5.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0423'
$ws.Range('E39').Value = '  +4.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.00'
$ws.Range('E40').Value = '  +1.73%  '
$ws.Range('D41').Value = '3.113.11'
$ws.Range('E41').Value = '  +4.99%  '
$ws.Range('E42').Value = '  +5.47%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.298'
$ws.Range('E43').Value = '  +15.03%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.121'
$ws.Range('E44').Value = '  +6.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.49'
$ws.Range('E45').Value = '  +17.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.17'
$ws.Range('E46').Value = '  +5.08%  '
$ws.Range('D47').Value = '0.0₃0578'
$ws.Range('E47').Value = '  +13.42%  '
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('E49').Value = '  +3.78%  '
$ws.Range('E50').Value = '  +11.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.93'
$ws.Range('E51').Value = '  -0.25%  '
